# Updated CHE_grids model - 2025-08-13 18:21
# Re-applies the refreshed "existing_stock" values (comm-out / ncap_pasti / ncap_cost /
# ncap_fom / act_cost) for rows 14-186, matching the latest VerveStacks CHE export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Cells.Item(14,5).Value = 0.066
$ws.Cells.Item(14,8).Value = 60.500000000000014
$ws.Cells.Item(15,5).Value = 0.086
$ws.Cells.Item(15,8).Value = 60.50000000000001
$ws.Cells.Item(16,3).Value = 'e_CH21-220'
$ws.Cells.Item(18,3).Value = 'e_w212722603-220'
$ws.Cells.Item(20,3).Value = 'e_CH21-220'
$ws.Cells.Item(22,3).Value = 'e_CH21-220'
$ws.Cells.Item(22,5).Value = 0.03992778473916777
$ws.Cells.Item(22,7).Value = 3712.5
$ws.Cells.Item(22,8).Value = 71.5
$ws.Cells.Item(22,9).Value = 2.52
$ws.Cells.Item(23,3).Value = 'e_CH50-220'
$ws.Cells.Item(23,5).Value = 0.018706720306008073
$ws.Cells.Item(23,7).Value = 3712.5
$ws.Cells.Item(23,8).Value = 71.5
$ws.Cells.Item(23,9).Value = 2.52
$ws.Cells.Item(24,3).Value = 'e_CH12-220'
$ws.Cells.Item(24,5).Value = 0.025445162566774423
$ws.Cells.Item(24,7).Value = 3712.5000000000005
$ws.Cells.Item(24,8).Value = 71.5
$ws.Cells.Item(24,9).Value = 2.52
$ws.Cells.Item(25,3).Value = 'e_w212722603-220'
$ws.Cells.Item(25,5).Value = 0.09413704412055675
$ws.Cells.Item(25,7).Value = 3162.5000000000005
$ws.Cells.Item(25,8).Value = 60.500000000000014
$ws.Cells.Item(25,9).Value = 2.3100000000000005
$ws.Cells.Item(26,3).Value = 'e_CH46-220'
$ws.Cells.Item(26,5).Value = 0.015790081118512187
$ws.Cells.Item(27,3).Value = 'e_CH30-380'
$ws.Cells.Item(27,5).Value = 0.013577458286618762
$ws.Cells.Item(27,7).Value = 3712.5000000000005
$ws.Cells.Item(28,3).Value = 'e_r7933294-380'
$ws.Cells.Item(28,5).Value = 0.07211138956670854
$ws.Cells.Item(28,7).Value = 3162.5000000000005
$ws.Cells.Item(28,8).Value = 60.500000000000014
$ws.Cells.Item(28,9).Value = 2.3100000000000005
$ws.Cells.Item(29,3).Value = 'e_CH45-220'
$ws.Cells.Item(29,5).Value = 0.06034425905163894
$ws.Cells.Item(29,7).Value = 3162.5000000000005
$ws.Cells.Item(29,8).Value = 60.500000000000014
$ws.Cells.Item(29,9).Value = 2.3100000000000005
$ws.Cells.Item(30,3).Value = 'e_CH17-380'
$ws.Cells.Item(30,5).Value = 0.07160852074127821
$ws.Cells.Item(30,7).Value = 3162.5000000000005
$ws.Cells.Item(30,8).Value = 60.500000000000014
$ws.Cells.Item(30,9).Value = 2.3100000000000005
$ws.Cells.Item(31,3).Value = 'e_CH41-380'
$ws.Cells.Item(31,5).Value = 0.02645090021763507
$ws.Cells.Item(32,3).Value = 'e_CH1-220'
$ws.Cells.Item(32,5).Value = 0.015086064762909735
$ws.Cells.Item(32,7).Value = 3712.5000000000005
$ws.Cells.Item(32,8).Value = 71.5
$ws.Cells.Item(32,9).Value = 2.52
$ws.Cells.Item(33,5).Value = 0.03821803073270466
$ws.Cells.Item(34,3).Value = 'e_w758943072-220'
$ws.Cells.Item(34,5).Value = 0.2540493306074
$ws.Cells.Item(35,3).Value = 'e_CH32-220'
$ws.Cells.Item(35,5).Value = 0.12340400976060163
$ws.Cells.Item(36,3).Value = 'e_CH19-220'
$ws.Cells.Item(36,5).Value = 0.30664940974741184
$ws.Cells.Item(36,7).Value = 2750.0
$ws.Cells.Item(36,8).Value = 55.00000000000001
$ws.Cells.Item(36,9).Value = 2.1
$ws.Cells.Item(37,3).Value = 'e_CH18-220'
$ws.Cells.Item(37,5).Value = 0.17932302314845372
$ws.Cells.Item(37,7).Value = 3162.5
$ws.Cells.Item(37,8).Value = 60.50000000000001
$ws.Cells.Item(38,3).Value = 'e_CH48-225'
$ws.Cells.Item(38,5).Value = 0.17017081052562183
$ws.Cells.Item(38,7).Value = 3162.5000000000005
$ws.Cells.Item(38,8).Value = 60.500000000000014
$ws.Cells.Item(38,9).Value = 2.3100000000000005
$ws.Cells.Item(40,3).Value = 'e_w212722603-220'
$ws.Cells.Item(42,3).Value = 'e_r7933294-380'
$ws.Cells.Item(45,3).Value = 'e_CH17-380'
$ws.Cells.Item(46,3).Value = 'e_CH50-220'
$ws.Cells.Item(50,3).Value = 'e_r7933294-380'
$ws.Cells.Item(58,3).Value = 'e_CH21-220'
$ws.Cells.Item(59,3).Value = 'e_w212722603-220'
$ws.Cells.Item(61,3).Value = 'e_r7933294-380'
$ws.Cells.Item(63,3).Value = 'e_CH45-220'
$ws.Cells.Item(65,3).Value = 'e_w212722603-220'
$ws.Cells.Item(66,3).Value = 'e_CH51-220'
$ws.Cells.Item(68,3).Value = 'e_w212722603-220'
$ws.Cells.Item(69,3).Value = 'e_CH17-380'
$ws.Cells.Item(72,3).Value = 'e_r7933294-380'
$ws.Cells.Item(73,3).Value = 'e_CH46-220'
$ws.Cells.Item(77,3).Value = 'e_CH21-220'
$ws.Cells.Item(84,3).Value = 'e_CH51-220'
$ws.Cells.Item(94,3).Value = 'e_CH46-220'
$ws.Cells.Item(98,3).Value = 'e_CH17-380'
$ws.Cells.Item(98,5).Value = 0.061
$ws.Cells.Item(98,7).Value = 2783.0
$ws.Cells.Item(98,8).Value = 66.55000000000001
$ws.Cells.Item(98,9).Value = 2.8875
$ws.Cells.Item(99,3).Value = 'e_CH1-220'
$ws.Cells.Item(99,5).Value = 0.05
$ws.Cells.Item(99,7).Value = 3267.0000000000005
$ws.Cells.Item(99,8).Value = 78.65
$ws.Cells.Item(99,9).Value = 3.1500000000000004
$ws.Cells.Item(100,3).Value = 'e_CH17-380'
$ws.Cells.Item(101,3).Value = 'e_CH50-220'
$ws.Cells.Item(104,3).Value = 'e_CH17-380'
$ws.Cells.Item(111,5).Value = 0.0011
$ws.Cells.Item(111,8).Value = 21.450000000000003
$ws.Cells.Item(113,5).Value = 0.0012
$ws.Cells.Item(114,5).Value = 0.0028
$ws.Cells.Item(114,8).Value = 21.450000000000006
$ws.Cells.Item(116,5).Value = 0.001
$ws.Cells.Item(117,3).Value = 'elc_spv-CHE_0022'
$ws.Cells.Item(117,5).Value = 0.0016
$ws.Cells.Item(118,3).Value = 'elc_spv-CHE_0013'
$ws.Cells.Item(118,5).Value = 0.0015
$ws.Cells.Item(119,3).Value = 'elc_spv-CHE_0003'
$ws.Cells.Item(119,5).Value = 0.0027
$ws.Cells.Item(119,7).Value = 1336.5
$ws.Cells.Item(120,3).Value = 'elc_spv-CHE_0006'
$ws.Cells.Item(120,5).Value = 0.001
$ws.Cells.Item(121,3).Value = 'elc_spv-CHE_0000'
$ws.Cells.Item(121,5).Value = 0.0018
$ws.Cells.Item(121,7).Value = 1336.5000000000002
$ws.Cells.Item(122,3).Value = 'elc_spv-CHE_0021'
$ws.Cells.Item(122,5).Value = 0.0017
$ws.Cells.Item(122,7).Value = 1336.4999999999998
$ws.Cells.Item(123,3).Value = 'elc_spv-CHE_0007'
$ws.Cells.Item(123,5).Value = 0.0025999999999999994
$ws.Cells.Item(124,3).Value = 'elc_spv-CHE_0004'
$ws.Cells.Item(124,5).Value = 0.0045
$ws.Cells.Item(125,3).Value = 'elc_spv-CHE_0008'
$ws.Cells.Item(125,5).Value = 0.002
$ws.Cells.Item(126,3).Value = 'elc_spv-CHE_0009'
$ws.Cells.Item(126,5).Value = 0.0017
$ws.Cells.Item(126,7).Value = 1336.4999999999998
$ws.Cells.Item(127,3).Value = 'elc_spv-CHE_0010'
$ws.Cells.Item(127,5).Value = 0.0058000000000000005
$ws.Cells.Item(128,3).Value = 'elc_spv-CHE_0012'
$ws.Cells.Item(128,5).Value = 0.0012
$ws.Cells.Item(129,3).Value = 'elc_spv-CHE_0014'
$ws.Cells.Item(129,5).Value = 0.0015
$ws.Cells.Item(129,7).Value = 1336.5
$ws.Cells.Item(130,5).Value = 0.0012
$ws.Cells.Item(131,3).Value = 'elc_spv-CHE_0006'
$ws.Cells.Item(131,5).Value = 0.0012
$ws.Cells.Item(132,3).Value = 'elc_spv-CHE_0004'
$ws.Cells.Item(132,5).Value = 0.024400000000000005
$ws.Cells.Item(132,7).Value = 1336.4999999999998
$ws.Cells.Item(132,8).Value = 21.45
$ws.Cells.Item(133,3).Value = 'elc_spv-CHE_0002'
$ws.Cells.Item(133,5).Value = 0.0012
$ws.Cells.Item(134,3).Value = 'elc_spv-CHE_0012'
$ws.Cells.Item(134,5).Value = 0.0013
$ws.Cells.Item(134,7).Value = 1336.5
$ws.Cells.Item(134,8).Value = 21.450000000000003
$ws.Cells.Item(135,3).Value = 'elc_spv-CHE_0007'
$ws.Cells.Item(135,5).Value = 0.0011
$ws.Cells.Item(136,3).Value = 'elc_spv-CHE_0022'
$ws.Cells.Item(136,5).Value = 0.001
$ws.Cells.Item(137,3).Value = 'elc_spv-CHE_0003'
$ws.Cells.Item(137,5).Value = 0.001
$ws.Cells.Item(138,3).Value = 'elc_spv-CHE_0004'
$ws.Cells.Item(139,3).Value = 'elc_spv-CHE_0003'
$ws.Cells.Item(139,5).Value = 0.0024000000000000002
$ws.Cells.Item(140,3).Value = 'elc_spv-CHE_0018'
$ws.Cells.Item(140,5).Value = 0.0013
$ws.Cells.Item(141,5).Value = 0.0038
$ws.Cells.Item(141,8).Value = 21.450000000000006
$ws.Cells.Item(142,3).Value = 'elc_spv-CHE_0006'
$ws.Cells.Item(142,5).Value = 0.001
$ws.Cells.Item(142,8).Value = 21.450000000000003
$ws.Cells.Item(143,3).Value = 'elc_spv-CHE_0001'
$ws.Cells.Item(143,5).Value = 0.0021000000000000003
$ws.Cells.Item(143,8).Value = 21.450000000000003
$ws.Cells.Item(144,3).Value = 'elc_spv-CHE_0012'
$ws.Cells.Item(144,5).Value = 0.005
$ws.Cells.Item(144,8).Value = 21.450000000000003
$ws.Cells.Item(145,5).Value = 0.0079
$ws.Cells.Item(145,8).Value = 21.45
$ws.Cells.Item(146,3).Value = 'elc_spv-CHE_0021'
$ws.Cells.Item(146,5).Value = 0.0079
$ws.Cells.Item(146,8).Value = 21.45
$ws.Cells.Item(147,3).Value = 'elc_spv-CHE_0022'
$ws.Cells.Item(148,3).Value = 'elc_spv-CHE_0021'
$ws.Cells.Item(148,5).Value = 0.004900000000000001
$ws.Cells.Item(149,3).Value = 'elc_spv-CHE_0014'
$ws.Cells.Item(149,5).Value = 0.0011
$ws.Cells.Item(150,3).Value = 'elc_spv-CHE_0002'
$ws.Cells.Item(150,5).Value = 0.0108
$ws.Cells.Item(150,8).Value = 21.45
$ws.Cells.Item(151,3).Value = 'elc_spv-CHE_0013'
$ws.Cells.Item(152,3).Value = 'elc_spv-CHE_0018'
$ws.Cells.Item(152,5).Value = 0.0014
$ws.Cells.Item(152,8).Value = 21.450000000000003
$ws.Cells.Item(153,3).Value = 'elc_spv-CHE_0009'
$ws.Cells.Item(153,5).Value = 0.004900000000000001
$ws.Cells.Item(154,5).Value = 0.0108
$ws.Cells.Item(154,8).Value = 21.45
$ws.Cells.Item(155,3).Value = 'elc_spv-CHE_0004'
$ws.Cells.Item(155,5).Value = 0.0011
$ws.Cells.Item(155,8).Value = 21.450000000000003
$ws.Cells.Item(156,3).Value = 'elc_spv-CHE_0008'
$ws.Cells.Item(157,3).Value = 'elc_spv-CHE_0002'
$ws.Cells.Item(157,5).Value = 0.20916506699980644
$ws.Cells.Item(158,3).Value = 'elc_spv-CHE_0024'
$ws.Cells.Item(158,5).Value = 0.1986733343920766
$ws.Cells.Item(159,3).Value = 'elc_spv-CHE_0025'
$ws.Cells.Item(159,5).Value = 0.1637027705855386
$ws.Cells.Item(160,3).Value = 'elc_spv-CHE_0006'
$ws.Cells.Item(160,5).Value = 0.2084203530758958
$ws.Cells.Item(161,3).Value = 'elc_spv-CHE_0011'
$ws.Cells.Item(161,5).Value = 0.15348244065625796
$ws.Cells.Item(162,3).Value = 'elc_spv-CHE_0000'
$ws.Cells.Item(162,5).Value = 0.21232141441249713
$ws.Cells.Item(163,3).Value = 'elc_spv-CHE_0019'
$ws.Cells.Item(163,5).Value = 0.16452449149779813
$ws.Cells.Item(164,3).Value = 'elc_spv-CHE_0005'
$ws.Cells.Item(164,5).Value = 0.135571742097842
$ws.Cells.Item(165,3).Value = 'elc_spv-CHE_0009'
$ws.Cells.Item(165,5).Value = 0.2148926965356589
$ws.Cells.Item(166,5).Value = 0.13455092994741213
$ws.Cells.Item(167,3).Value = 'elc_spv-CHE_0023'
$ws.Cells.Item(167,5).Value = 0.156849148622487
$ws.Cells.Item(168,5).Value = 0.15237710864056714
$ws.Cells.Item(169,3).Value = 'elc_spv-CHE_0012'
$ws.Cells.Item(169,5).Value = 0.1956202890369272
$ws.Cells.Item(170,3).Value = 'elc_spv-CHE_0001'
$ws.Cells.Item(170,5).Value = 0.1516718373204449
$ws.Cells.Item(171,3).Value = 'elc_spv-CHE_0018'
$ws.Cells.Item(171,5).Value = 0.19607575719972248
$ws.Cells.Item(172,3).Value = 'elc_spv-CHE_0020'
$ws.Cells.Item(172,5).Value = 0.1609643269510828
$ws.Cells.Item(173,3).Value = 'elc_spv-CHE_0021'
$ws.Cells.Item(173,5).Value = 0.19094540137838087
$ws.Cells.Item(174,3).Value = 'elc_spv-CHE_0015'
$ws.Cells.Item(174,5).Value = 0.16891764627214115
$ws.Cells.Item(175,3).Value = 'elc_spv-CHE_0008'
$ws.Cells.Item(175,5).Value = 0.15972751958943046
$ws.Cells.Item(176,3).Value = 'elc_spv-CHE_0022'
$ws.Cells.Item(176,5).Value = 0.19396276056819794
$ws.Cells.Item(177,3).Value = 'elc_spv-CHE_0004'
$ws.Cells.Item(177,5).Value = 0.18104249049863588
$ws.Cells.Item(178,3).Value = 'elc_spv-CHE_0010'
$ws.Cells.Item(178,5).Value = 0.19113510152041108
$ws.Cells.Item(179,3).Value = 'elc_spv-CHE_0007'
$ws.Cells.Item(179,5).Value = 0.16513303452767722
$ws.Cells.Item(180,3).Value = 'elc_spv-CHE_0017'
$ws.Cells.Item(180,5).Value = 0.15120603895189552
$ws.Cells.Item(181,3).Value = 'elc_spv-CHE_0014'
$ws.Cells.Item(181,5).Value = 0.17086629872121376
$ws.Cells.Item(185,3).Value = 'elc_won-CHE_0006'
$ws.Cells.Item(186,3).Value = 'elc_won-CHE_0010'
